$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A12's slightly-adjusted timestamp value
$ws.Range("A12").Value = 45878.50017998843

# Add new row 13 data
$ws.Range("A13").Value = 45878.54183790438
$ws.Range("B13").Value = 2025
$ws.Range("C13").Value = 37
$ws.Range("D13").Value = 17.33
$ws.Range("E13").Value = 82.66
$ws.Range("F13").Value = 615.23
$ws.Range("G13").Value = 13.45
$ws.Range("H13").Value = "ESE"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "13:00:14"

# Match the date/time number format used by the other rows in column A
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
